# The model was not picking up the correct soil depth value from the test
# configs: the "Texture" label used the non-matching value "Silt loam"
# (should be "Silt") and the "SampleDepth" label used "0-30cm" (should be
# "Top30cm"). Update the FieldConfigs test sheet accordingly for every
# scenario column (C:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 = Texture -> "Silt" (was "Silt loam")
$ws.Range("C5:L5").Value = "Silt"

# Row 7 = SampleDepth -> "Top30cm" (was "0-30cm")
$ws.Range("C7:L7").Value = "Top30cm"

# Reflect the edit location as the active selection, matching where the
# author was last working when the workbook was saved.
$ws.Range("C7:L7").Select() | Out-Null
